# ---------------------------------------------------------------------
# New 2022-Q4 fund-holdings data (fund code / name / size / position /
# weight / market value are kept as text, exactly like the existing
# quarter sheets, except the two zero market-value rows which store a
# real number 0).
# ---------------------------------------------------------------------
$q4Data = @(
    @(0,  "002938", "中银证券健康产业灵活配置混合",       "5.33", "92.95", "4.03", "0.2148", 7),
    @(1,  "002601", "中银证券价值精选灵活配置混合",       "5.15", "92.90", "3.64", "0.1875", 10),
    @(2,  "519956", "长信睿进灵活配置混合C",             "6.44", "43.59", "2.26", "0.1455", 9),
    @(3,  "011346", "淳厚鑫淳一年持有期混合",             "3.34", "78.96", "3.13", "0.1045", 3),
    @(4,  "501038", "银华明择多策略定期开放混合",         "3.11", "81.75", "3.07", "0.0955", 9),
    @(5,  "519120", "浦银安盛新兴产业混合A",             "1.72", "91.65", "2.78", "0.0478", 6),
    @(6,  "012454", "淳厚鑫悦混合A",                     "1.82", "85.29", "2.57", "0.0468", 7),
    @(7,  "519113", "浦银安盛精致生活混合",               "1.81", "90.69", "2.37", "0.0429", 10),
    @(8,  "011270", "中银证券优势制造股票C",             "1.09", "93.20", "3.80", "0.0414", 8),
    @(9,  "011269", "中银证券优势制造股票A",             "0.82", "93.20", "3.80", "0.0312", 8),
    @(10, "012455", "淳厚鑫悦混合C",                     "0.57", "85.29", "2.57", "0.0146", 7),
    @(11, "519175", "浦银安盛经济带崛起灵活配置混合",     "0.39", "23.04", "1.05", "0.0041", 5),
    @(12, "004801", "浦银安盛安久回报定期开放混合A",     "0.17", "25.01", "1.29", "0.0022", 7),
    @(13, "008162", "浦银安盛经济带崛起灵活配置混合C",   "0.12", "23.04", "1.05", "0.0013", 5),
    @(14, "519957", "长信睿进灵活配置混合A",             "0.01", "43.59", "2.26", "0.0002", 9),
    @(15, "004802", "浦银安盛安久回报定期开放混合C",     "0.00", "25.01", "1.29", "0.00",   7),
    @(16, "014061", "浦银安盛新兴产业混合C",             "0.00", "91.65", "2.78", "0.00",   6)
)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "总计" (summary) sheet: insert a new row 2 for 2022-Q4 and shift
#    the existing quarters down by one row.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")
$summary.Rows.Item(2).Insert()

# Restore the formatting "Insert" didn't carry onto the fresh row 2 by
# copying the format of the row directly below (still the pristine,
# untouched original formatting of the old row 2).
$summary.Range("A3:D3").Copy()
$summary.Range("A2:D2").PasteSpecial(-4122)

$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q4"
$summary.Cells.Item(2, 3).Value = 17
$summary.Cells.Item(2, 4).Value = 0.98

# Column A is a plain 0-based row index (row2 -> 0, row3 -> 1, ...); it
# doesn't travel with the shifted rows, so renumber rows 3-6.
$summary.Cells.Item(3, 1).Value = 1
$summary.Cells.Item(4, 1).Value = 2
$summary.Cells.Item(5, 1).Value = 3
$summary.Cells.Item(6, 1).Value = 4

# ---------------------------------------------------------------------
# 2. Insert a brand-new worksheet "2022-Q4" before "2022-Q3" by copying
#    "2022-Q3" (so the header/number formatting carries over) and then
#    overwriting its data with the 2022-Q4 figures.
# ---------------------------------------------------------------------
$sheetQ3 = $wb.Worksheets.Item("2022-Q3")
$sheetQ3.Copy($sheetQ3, $null)
$q4 = $wb.Worksheets.Item("2022-Q3 (2)")
$q4.Name = "2022-Q4"

# The Q3 sheet has 22 data rows (rows 2-23); Q4 only has 17 (rows 2-18).
# Drop the now-unused trailing rows.
$q4.Range("A19:H23").EntireRow.Delete()

# Fund code/name/size/position/weight/market-value columns hold text
# (even though many look numeric) so leading zeros and trailing zeros
# in codes/percentages survive - force text format before writing.
$q4.Range("B2:G18").NumberFormat = "@"

for ($i = 0; $i -lt $q4Data.Length; $i++) {
    $r = $i + 2
    $row = $q4Data[$i]

    $q4.Cells.Item($r, 1).Value = $row[0]
    $q4.Cells.Item($r, 2).Value = $row[1]
    $q4.Cells.Item($r, 3).Value = $row[2]
    $q4.Cells.Item($r, 4).Value = $row[3]
    $q4.Cells.Item($r, 5).Value = $row[4]
    $q4.Cells.Item($r, 6).Value = $row[5]
    $q4.Cells.Item($r, 7).Value = $row[6]
    $q4.Cells.Item($r, 8).Value = $row[7]
}

# The last two rows' market value is exactly 0 - stored as a real
# number rather than the text "0.00".
$q4.Range("G17:G18").NumberFormat = "General"
$q4.Cells.Item(17, 7).Value = 0
$q4.Cells.Item(18, 7).Value = 0
